$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-G)
$data = @(
    @{Row=2;  A=83444; B="Ana Clara Silveira";    C="TI";         D="Doenca";              E=4; F=45095; G=7827.71},
    @{Row=3;  A=23902; B="Maria Alice Rios";       C="Juridico";   D="Problemas pessoais";  E=1; F=45079; G=3438.82},
    @{Row=4;  A=6210;  B="Maya Pastor";             C="Vendas";     D="Outros";              E=2; F=45092; G=2426.54},
    @{Row=5;  A=17974; B="Lavínia da Mata";         C="Financeiro"; D="Viagem de negocios";  E=1; F=45089; G=7865.29},
    @{Row=6;  A=7165;  B="Sr. Ravi Barbosa";        C="TI";         D="Problemas pessoais";  E=4; F=45092; G=2143.8},
    @{Row=7;  A=28089; B="Clarice Caldeira";        C="TI";         D="Consulta medica";     E=7; F=45094; G=2407.57},
    @{Row=8;  A=63028; B="Giovanna Pacheco";        C="Marketing";  D="Viagem de negocios";  E=5; F=45099; G=9355.99},
    @{Row=9;  A=49078; B="Maria Cecília da Paz";    C="TI";         D="Problemas pessoais";  E=2; F=45079; G=8634.9},
    @{Row=10; A=54664; B="José Pedro Martins";      C="Financeiro"; D="Outros";              E=3; F=45083; G=9292.76},
    @{Row=11; A=75701; B="Melina da Mata";          C="Operacoes";  D="Outros";              E=8; F=45098; G=9178.83}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
